$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text37 = "LOT2052 -  Tecnologia de Bebidas Experimental  (Indicação de Conjunto)`n"
$text38 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"

# Swap the shared-string contents referenced by row 23 (was LOT2052) and row 24 (was LOT2028)
# so that row 23 now shows the LOT2028 text and row 24 shows the LOT2052 text, matching
# the reordering of <si> entries in the diff.
$ws.Range("B23").Value = $text38
$ws.Range("C23").Value = $text38
$ws.Range("B24").Value = $text37
$ws.Range("C24").Value = $text37
